$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "I have a two math degrees." -> "I have two math degrees."
# -----------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("I have a two math degrees.", $true, $false, $false, $false, $false,
                    $true, 1, $false, "I have two math degrees.", 2)

# -----------------------------------------------------------------
# Change 2: "October 1st." -> "October 1" + superscript "st" + "."
# -----------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("1st", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found2) {
    # Narrow the found range down to just the "st" portion (skip the "1").
    $stRange = $d.Range($rng2.Start + 1, $rng2.End)
    $stRange.Font.Superscript = $true
}

# -----------------------------------------------------------------
# Change 3: collapse the split runs "know" + "s" and "w" + "a" + "s"
# back into plain merged text (no formatting change), i.e. the
# paragraph's text ends up identical but stored as a single run.
# -----------------------------------------------------------------
$rng3 = $d.Content
$finalText = "He claims spiraling debt; funny, it rose 8.3% per year under Biden, 8.9% per year under Trump.  He blames Biden for the invasion of Ukraine, when everyone knows that Trump's puppet-master, Putin, was responsible.  He mentions a Biden gaffe; the continuing stream of gaffes from Trump is infinitely worse."
$found3 = $rng3.Find.Execute($finalText, $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found3) {
    $start3 = $rng3.Start
    $end3 = $rng3.End
    # Overwrite with a placeholder first, then restore the real text, so
    # the run-rebuild is forced to happen and all sub-runs coalesce into one.
    $rng3.Text = "X"
    $mergedRange = $d.Range($start3, $start3 + 1)
    $mergedRange.Text = $finalText
}
